# Add a new "2022-Q4" quarterly sheet (copied from "2022-Q3" to keep the same
# layout/formatting), positioned right after the "总计" (total) summary sheet,
# then fill it with the new quarter's fund-holding figures. Finally, record
# the new quarter in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q4" sheet by copying "2022-Q3" (same header row /
#    number formats / borders), inserted immediately before "2022-Q3".
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

# Columns D:G on this sheet hold numeric-looking figures stored as text;
# force text formatting so they keep their text type instead of turning
# into numbers when assigned below.
$wsQ4.Range("D2:G3").NumberFormat = "@"

# Fund 006923 - 前海开源沪港深非周期性行业股票A
$wsQ4.Range("D2").Value = "0.28"
$wsQ4.Range("E2").Value = "90.65"
$wsQ4.Range("F2").Value = "5.52"
$wsQ4.Range("G2").Value = "0.0155"
$wsQ4.Range("H2").Value = 7

# Fund 006924 - 前海开源沪港深非周期性行业股票C
$wsQ4.Range("D3").Value = "0.24"
$wsQ4.Range("E3").Value = "90.65"
$wsQ4.Range("F3").Value = "5.52"
$wsQ4.Range("G3").Value = "0.0132"
$wsQ4.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" (total) summary sheet: add a 2022-Q4 row at the
#    top of the data and push the rest of the quarters down by one row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Make room for one more data row (7 rows of data incl. header -> A1:D7).
$wsTotal.Rows.Item(6).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.03

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.05

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2022-Q1"
$wsTotal.Range("C5").Value = 2
$wsTotal.Range("D5").Value = 0.04

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q4"
$wsTotal.Range("C6").Value = 3
$wsTotal.Range("D6").Value = 0.16

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2021-Q3"
$wsTotal.Range("C7").Value = 4
$wsTotal.Range("D7").Value = 0.57

# The newly inserted row's index cell (A6) should look like the others in
# column A (bold, centred, thin border) - match that formatting; row 7's
# index cell already carried it over from the row that used to be there.
$wsTotal.Range("A6").Font.Bold = $true
$wsTotal.Range("A6").HorizontalAlignment = -4108
$wsTotal.Range("A6").VerticalAlignment = -4160
$wsTotal.Range("A6").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Keep "2021-Q3" (now the last tab) the selected/active sheet, same as
#    before the edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
